$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 104.46667
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H39").Value = 1340.6471
$ws.Range("I39").Value = 478.7
$ws.Range("J39").Value = 2572
$ws.Range("K39").Value = 1436.1
$ws.Range("L39").Value = 7716
$ws.Range("M39").Value = -1140.1
$ws.Range("N39").Value = -8308
$ws.Range("H42").Value = 2067.2727
$ws.Range("I42").Value = 420.42856
$ws.Range("K42").Value = 1261.28568
$ws.Range("M42").Value = -1031.28568
$ws.Range("H58").Value = 943.9231
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9300
$ws.Range("H82").Value = 5761.4287
$ws.Range("I82").Value = 3221.8333
$ws.Range("K82").Value = 9665.499899999999
$ws.Range("M82").Value = -9259.499899999999
$ws.Range("H85").Value = 5761.4287
$ws.Range("I85").Value = 3221.8333
$ws.Range("K85").Value = 9665.499899999999
$ws.Range("M85").Value = -8261.499899999999
$ws.Range("H98").Value = 2609.3845
$ws.Range("I98").Value = 2083.9092
$ws.Range("K98").Value = 2083.9092
$ws.Range("M98").Value = -585.9092000000001
$ws.Range("H100").Value = 2987.5
$ws.Range("I100").Value = 2294.25
$ws.Range("K100").Value = 2294.25
$ws.Range("M100").Value = -1753.25
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -36988
$ws.Range("H111").Value = 665.4286
$ws.Range("I111").Value = 609.6667
$ws.Range("K111").Value = 1829.0001
$ws.Range("M111").Value = 1237.9999
$ws.Range("H112").Value = 2555.25
$ws.Range("J112").Value = 2714.9092
$ws.Range("L112").Value = 8144.7276
$ws.Range("N112").Value = -10360.7276
$ws.Range("H116").Value = 4725.952
$ws.Range("I116").Value = 5795.923
$ws.Range("J116").Value = 2987.25
$ws.Range("K116").Value = 5795.923
$ws.Range("L116").Value = 2987.25
$ws.Range("M116").Value = -2353.923
$ws.Range("N116").Value = -9871.25
$ws.Range("H122").Value = 2609.3845
$ws.Range("I122").Value = 2083.9092
$ws.Range("K122").Value = 6251.7276
$ws.Range("M122").Value = -3801.7276
$ws.Range("H125").Value = 4500.2354
$ws.Range("I125").Value = 3247
$ws.Range("J125").Value = 4885.846
$ws.Range("K125").Value = 29223
$ws.Range("L125").Value = 43972.61399999999
$ws.Range("M125").Value = -26763
$ws.Range("N125").Value = -48892.61399999999
$ws.Range("H132").Value = 5109.5
$ws.Range("I132").Value = 4807.5
$ws.Range("J132").Value = 6116.1665
$ws.Range("K132").Value = 14422.5
$ws.Range("L132").Value = 18348.4995
$ws.Range("M132").Value = -11892.5
$ws.Range("N132").Value = -23408.4995
$ws.Range("H137").Value = 29356.334
$ws.Range("I137").Value = 50209.57
$ws.Range("J137").Value = 5027.5557
$ws.Range("K137").Value = 150628.71
$ws.Range("L137").Value = 15082.6671
$ws.Range("M137").Value = -148078.71
$ws.Range("N137").Value = -20182.6671
$ws.Range("H138").Value = 3961.9707
$ws.Range("J138").Value = 3970.7666
$ws.Range("L138").Value = 11912.2998
$ws.Range("N138").Value = -22192.2998
$ws.Range("H141").Value = 1171.862
$ws.Range("I141").Value = 1172.5
$ws.Range("J141").Value = 1166.3334
$ws.Range("K141").Value = 3517.5
$ws.Range("L141").Value = 3499.0002
$ws.Range("M141").Value = 1662.5
$ws.Range("N141").Value = -13859.0002

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1638.6666
$ws.Range("I2").Value = 1321.6316
$ws.Range("J2").Value = 2391.625
$ws.Range("K2").Value = 1321.6316
$ws.Range("L2").Value = 2391.625
$ws.Range("M2").Value = -1208.6316
$ws.Range("N2").Value = -2617.625
$ws.Range("H45").Value = 5865.778
$ws.Range("I45").Value = 3774.25
$ws.Range("K45").Value = 3774.25
$ws.Range("M45").Value = -3397.25
$ws.Range("H61").Value = 3458.0476
$ws.Range("I61").Value = 3357.5
$ws.Range("J61").Value = 3779.8
$ws.Range("K61").Value = 3357.5
$ws.Range("L61").Value = 3779.8
$ws.Range("M61").Value = -3145.5
$ws.Range("N61").Value = -4203.8
$ws.Range("H63").Value = 2197.3333
$ws.Range("J63").Value = 2196
$ws.Range("L63").Value = 2196
$ws.Range("N63").Value = -3568
$ws.Range("H66").Value = 2197.3333
$ws.Range("J66").Value = 2196
$ws.Range("L66").Value = 10980
$ws.Range("N66").Value = -17844
$ws.Range("H68").Value = 112545
$ws.Range("I68").Value = 75090
$ws.Range("K68").Value = 75090
$ws.Range("M68").Value = -74279
$ws.Range("H71").Value = 112545
$ws.Range("I71").Value = 75090
$ws.Range("K71").Value = 225270
$ws.Range("M71").Value = -221214
$ws.Range("H74").Value = 92443.37
$ws.Range("I74").Value = 101637.8
$ws.Range("J74").Value = 499
$ws.Range("K74").Value = 101637.8
$ws.Range("L74").Value = 499
$ws.Range("M74").Value = -100763.8
$ws.Range("N74").Value = -2247
$ws.Range("H77").Value = 92443.37
$ws.Range("I77").Value = 101637.8
$ws.Range("J77").Value = 499
$ws.Range("K77").Value = 508189
$ws.Range("L77").Value = 2495
$ws.Range("M77").Value = -503821
$ws.Range("N77").Value = -11231
$ws.Range("H88").Value = 3748.0908
$ws.Range("I88").Value = 2634.3333
$ws.Range("J88").Value = 4165.75
$ws.Range("K88").Value = 2634.3333
$ws.Range("L88").Value = 4165.75
$ws.Range("M88").Value = -2228.3333
$ws.Range("N88").Value = -4977.75
$ws.Range("H91").Value = 3748.0908
$ws.Range("I91").Value = 2634.3333
$ws.Range("J91").Value = 4165.75
$ws.Range("K91").Value = 2634.3333
$ws.Range("L91").Value = 4165.75
$ws.Range("M91").Value = -1230.3333
$ws.Range("N91").Value = -6973.75
$ws.Range("H102").Value = 3086.3333
$ws.Range("I102").Value = 2129.5
$ws.Range("K102").Value = 2129.5
$ws.Range("M102").Value = -507.5
$ws.Range("H116").Value = 1638.6666
$ws.Range("I116").Value = 1321.6316
$ws.Range("J116").Value = 2391.625
$ws.Range("K116").Value = 1321.6316
$ws.Range("L116").Value = 2391.625
$ws.Range("M116").Value = 972.3684000000001
$ws.Range("N116").Value = -6979.625
$ws.Range("H132").Value = 19617.105
$ws.Range("I132").Value = 21712.705
$ws.Range("J132").Value = 1804.5
$ws.Range("K132").Value = 65138.11500000001
$ws.Range("L132").Value = 5413.5
$ws.Range("M132").Value = -62608.11500000001
$ws.Range("N132").Value = -10473.5
$ws.Range("H136").Value = 3458.0476
$ws.Range("I136").Value = 3357.5
$ws.Range("J136").Value = 3779.8
$ws.Range("K136").Value = 10072.5
$ws.Range("L136").Value = 11339.4
$ws.Range("M136").Value = -7522.5
$ws.Range("N136").Value = -16439.4
$ws.Range("H137").Value = 99999
$ws.Range("J137").Value = 99999
$ws.Range("L137").Value = 99999
$ws.Range("N137").Value = -110199

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1638.6666
$ws.Range("I3").Value = 1321.6316
$ws.Range("J3").Value = 2391.625
$ws.Range("K3").Value = 1321.6316
$ws.Range("L3").Value = 2391.625
$ws.Range("M3").Value = -1207.6316
$ws.Range("N3").Value = -2619.625
$ws.Range("H8").Value = 300
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H10").Value = 9999
$ws.Range("J10").Value = 9999
$ws.Range("L10").Value = 9999
$ws.Range("N10").Value = -10279
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H59").Value = 99999
$ws.Range("J59").Value = 99999
$ws.Range("L59").Value = 99999
$ws.Range("N59").Value = -101693
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H99").Value = 86400.25
$ws.Range("I99").Value = 114200.336
$ws.Range("K99").Value = 114200.336
$ws.Range("M99").Value = -112702.336
$ws.Range("H108").Value = 45000
$ws.Range("J108").Value = 45000
$ws.Range("L108").Value = 45000
$ws.Range("N108").Value = -52680
$ws.Range("H123").Value = 59999
$ws.Range("J123").Value = 59999
$ws.Range("L123").Value = 59999
$ws.Range("N123").Value = -69799
$ws.Range("H134").Value = 3518.611
$ws.Range("I134").Value = 3378.7812
$ws.Range("K134").Value = 10136.3436
$ws.Range("M134").Value = -7601.3436

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1256.2667
$ws.Range("I16").Value = 804.8889
$ws.Range("K16").Value = 804.8889
$ws.Range("M16").Value = -517.8889
$ws.Range("H28").Value = 86666.664
$ws.Range("J28").Value = 86666.664
$ws.Range("L28").Value = 86666.664
$ws.Range("N28").Value = -87156.664
$ws.Range("H31").Value = 2667.35
$ws.Range("I31").Value = 2138.647
$ws.Range("K31").Value = 2138.647
$ws.Range("M31").Value = -1843.647
$ws.Range("H32").Value = 2577.25
$ws.Range("I32").Value = 1769.6666
$ws.Range("K32").Value = 1769.6666
$ws.Range("M32").Value = -1453.6666
$ws.Range("H34").Value = 2667.35
$ws.Range("I34").Value = 2138.647
$ws.Range("K34").Value = 2138.647
$ws.Range("M34").Value = -1936.647
$ws.Range("H58").Value = 50130.715
$ws.Range("I58").Value = 73750.92999999999
$ws.Range("K58").Value = 73750.92999999999
$ws.Range("M58").Value = -73547.92999999999
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 5000
$ws.Range("K86").Value = 5000
$ws.Range("M86").Value = -3877
$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 5000
$ws.Range("K89").Value = 25000
$ws.Range("M89").Value = -19384
$ws.Range("H105").Value = 3037.5
$ws.Range("I105").Value = 3037.5
$ws.Range("K105").Value = 3037.5
$ws.Range("M105").Value = -1290.5
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H107").Value = 445.8
$ws.Range("I107").Value = 461.7857
$ws.Range("K107").Value = 461.7857
$ws.Range("M107").Value = 1458.2143
$ws.Range("H113").Value = 1256.2667
$ws.Range("I113").Value = 804.8889
$ws.Range("K113").Value = 804.8889
$ws.Range("M113").Value = 1365.1111
$ws.Range("H132").Value = 2026.8422
$ws.Range("I132").Value = 1795.7858
$ws.Range("J132").Value = 2673.8
$ws.Range("K132").Value = 5387.357400000001
$ws.Range("L132").Value = 8021.400000000001
$ws.Range("M132").Value = -2857.357400000001
$ws.Range("N132").Value = -13081.4
$ws.Range("H136").Value = 50130.715
$ws.Range("I136").Value = 73750.92999999999
$ws.Range("K136").Value = 221252.79
$ws.Range("M136").Value = -218702.79

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 1058.4166
$ws.Range("I6").Value = 1220.9
$ws.Range("J6").Value = 246
$ws.Range("K6").Value = 3662.7
$ws.Range("L6").Value = 738
$ws.Range("M6").Value = -3549.7
$ws.Range("N6").Value = -964
$ws.Range("H7").Value = 3032.25
$ws.Range("I7").Value = 19
$ws.Range("J7").Value = 4036.6667
$ws.Range("K7").Value = 57
$ws.Range("L7").Value = 12110.0001
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = -12334.0001
$ws.Range("H70").Value = 16025
$ws.Range("I70").Value = 12100
$ws.Range("J70").Value = 19950
$ws.Range("K70").Value = 36300
$ws.Range("L70").Value = 59850
$ws.Range("M70").Value = -35985
$ws.Range("N70").Value = -60480
$ws.Range("H73").Value = 16025
$ws.Range("I73").Value = 12100
$ws.Range("J73").Value = 19950
$ws.Range("K73").Value = 36300
$ws.Range("L73").Value = 59850
$ws.Range("M73").Value = -35208
$ws.Range("N73").Value = -62034
$ws.Range("H92").Value = 318.6154
$ws.Range("I92").Value = 299.5
$ws.Range("J92").Value = 382.33334
$ws.Range("K92").Value = 898.5
$ws.Range("L92").Value = 1147.00002
$ws.Range("M92").Value = 349.5
$ws.Range("N92").Value = -3643.00002
$ws.Range("H93").Value = 10864
$ws.Range("I93").Value = 3024
$ws.Range("J93").Value = 14000
$ws.Range("K93").Value = 9072
$ws.Range("L93").Value = 42000
$ws.Range("M93").Value = -7200
$ws.Range("N93").Value = -45744
$ws.Range("H113").Value = 493.4375
$ws.Range("I113").Value = 476.625
$ws.Range("J113").Value = 510.25
$ws.Range("K113").Value = 1429.875
$ws.Range("L113").Value = 1530.75
$ws.Range("M113").Value = 740.125
$ws.Range("N113").Value = -5870.75
$ws.Range("H132").Value = 3449.75
$ws.Range("I132").Value = 1800
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 16200
$ws.Range("L132").Value = 35997.0003
$ws.Range("M132").Value = -13670
$ws.Range("N132").Value = -41057.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 119.947365
$ws.Range("I2").Value = 118.27778
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 118.27778
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = -5.277780000000007
$ws.Range("N2").Value = -376
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H70").Value = 7466.6665
$ws.Range("I70").Value = 6500
$ws.Range("J70").Value = 7950
$ws.Range("K70").Value = 6500
$ws.Range("L70").Value = 7950
$ws.Range("M70").Value = -6230
$ws.Range("N70").Value = -8490
$ws.Range("H73").Value = 7466.6665
$ws.Range("I73").Value = 6500
$ws.Range("J73").Value = 7950
$ws.Range("K73").Value = 6500
$ws.Range("L73").Value = 7950
$ws.Range("M73").Value = -5564
$ws.Range("N73").Value = -9822
$ws.Range("H97").Value = 883.9
$ws.Range("I97").Value = 675.2941
$ws.Range("J97").Value = 2066
$ws.Range("K97").Value = 675.2941
$ws.Range("L97").Value = 2066
$ws.Range("M97").Value = -179.2941
$ws.Range("N97").Value = -3058
$ws.Range("H102").Value = 3815
$ws.Range("I102").Value = 3910.6667
$ws.Range("J102").Value = 3599.75
$ws.Range("K102").Value = 3910.6667
$ws.Range("L102").Value = 3599.75
$ws.Range("M102").Value = -2288.6667
$ws.Range("N102").Value = -6843.75
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H113").Value = 103515.35
$ws.Range("I113").Value = 70187.39999999999
$ws.Range("K113").Value = 70187.39999999999
$ws.Range("M113").Value = -68017.39999999999
$ws.Range("H122").Value = 4066.9285
$ws.Range("I122").Value = 2448.182
$ws.Range("K122").Value = 7344.545999999999
$ws.Range("M122").Value = -4894.545999999999
$ws.Range("H126").Value = 7283.647
$ws.Range("I126").Value = 7070.3335
$ws.Range("K126").Value = 21211.0005
$ws.Range("M126").Value = -18741.0005
$ws.Range("H131").Value = 79500
$ws.Range("J131").Value = 79500
$ws.Range("L131").Value = 79500
$ws.Range("N131").Value = -89580
$ws.Range("H132").Value = 26491.701
$ws.Range("I132").Value = 35144.59
$ws.Range("J132").Value = 3861.077
$ws.Range("K132").Value = 105433.77
$ws.Range("L132").Value = 11583.231
$ws.Range("M132").Value = -102903.77
$ws.Range("N132").Value = -16643.231
$ws.Range("H137").Value = 70780
$ws.Range("J137").Value = 70780
$ws.Range("L137").Value = 70780
$ws.Range("N137").Value = -80980

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 8501
$ws.Range("I12").Value = 2002
$ws.Range("J12").Value = 15000
$ws.Range("K12").Value = 2002
$ws.Range("L12").Value = 15000
$ws.Range("M12").Value = -1832
$ws.Range("N12").Value = -15340
$ws.Range("H22").Value = 57627.4
$ws.Range("I22").Value = 185796.17
$ws.Range("K22").Value = 185796.17
$ws.Range("M22").Value = -185501.17
$ws.Range("H27").Value = 57627.4
$ws.Range("I27").Value = 185796.17
$ws.Range("K27").Value = 185796.17
$ws.Range("M27").Value = -185689.17
$ws.Range("H70").Value = 49999
$ws.Range("J70").Value = 49999
$ws.Range("L70").Value = 49999
$ws.Range("N70").Value = -50539
$ws.Range("H73").Value = 49999
$ws.Range("J73").Value = 49999
$ws.Range("L73").Value = 49999
$ws.Range("N73").Value = -51871
$ws.Range("H82").Value = 1887.0869
$ws.Range("I82").Value = 1143.6666
$ws.Range("J82").Value = 2365
$ws.Range("K82").Value = 1143.6666
$ws.Range("L82").Value = 2365
$ws.Range("M82").Value = -782.6666
$ws.Range("N82").Value = -3087
$ws.Range("H85").Value = 1887.0869
$ws.Range("I85").Value = 1143.6666
$ws.Range("J85").Value = 2365
$ws.Range("K85").Value = 1143.6666
$ws.Range("L85").Value = 2365
$ws.Range("M85").Value = 104.3334
$ws.Range("N85").Value = -4861
$ws.Range("H132").Value = 97423.62
$ws.Range("I132").Value = 124305.7
$ws.Range("J132").Value = 7816.6665
$ws.Range("K132").Value = 372917.1
$ws.Range("L132").Value = 23449.9995
$ws.Range("M132").Value = -370387.1
$ws.Range("N132").Value = -28509.9995
$ws.Range("H134").Value = 96598.39999999999
$ws.Range("J134").Value = 94999.25
$ws.Range("L134").Value = 94999.25
$ws.Range("N134").Value = -105139.25
$ws.Range("H137").Value = 99999
$ws.Range("J137").Value = 99999
$ws.Range("L137").Value = 99999
$ws.Range("N137").Value = -110199

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 50000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 50000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 50000
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -50228
$ws.Range("H4").Value = 16000
$ws.Range("I4").Value = 16000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 16000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -15887
$ws.Range("N4").ClearContents()
$ws.Range("H8").Value = 20000000
$ws.Range("I8").Value = 20000000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 20000000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -19999860
$ws.Range("N8").ClearContents()
$ws.Range("H46").Value = 56170.8
$ws.Range("J46").Value = 56170.8
$ws.Range("L46").Value = 56170.8
$ws.Range("N46").Value = -56632.8
$ws.Range("H96").Value = 5265.2856
$ws.Range("I96").Value = 3464.5
$ws.Range("J96").Value = 7666.3335
$ws.Range("K96").Value = 3464.5
$ws.Range("L96").Value = 7666.3335
$ws.Range("M96").Value = -2091.5
$ws.Range("N96").Value = -10412.3335
$ws.Range("H107").Value = 1531.15
$ws.Range("I107").Value = 965.0769
$ws.Range("J107").Value = 2582.4285
$ws.Range("K107").Value = 2895.2307
$ws.Range("L107").Value = 7747.2855
$ws.Range("M107").Value = -975.2307000000001
$ws.Range("N107").Value = -11587.2855
$ws.Range("H122").Value = 7500
$ws.Range("I122").Value = 10000
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 30000
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -27550
$ws.Range("N122").Value = -19900
$ws.Range("H126").Value = 74976.28999999999
$ws.Range("I126").Value = 102588.7
$ws.Range("K126").Value = 307766.1
$ws.Range("M126").Value = -305296.1
$ws.Range("H132").Value = 22578.736
$ws.Range("I132").Value = 25382.783
$ws.Range("J132").Value = 4152.143
$ws.Range("K132").Value = 76148.349
$ws.Range("L132").Value = 12456.429
$ws.Range("M132").Value = -73618.349
$ws.Range("N132").Value = -17516.429
$ws.Range("H134").Value = 56170.8
$ws.Range("J134").Value = 56170.8
$ws.Range("L134").Value = 168512.4
$ws.Range("N134").Value = -173582.4
